$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Deadly 5, a wild west-themed slot game featuring four outlaw characters. Play for free and experience the excitement of the American frontier.</w:t></w:r></w:p>'
$metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) Near the end of the document there used to be two paragraphs:
#       - bold:   "Play Deadly 5 Free - Wild West Themed Slot Review"
#       - italic: "Read our review of Deadly 5, ..."
#    The bold paragraph is removed entirely, and the italic
#    paragraph's text is replaced with the new image-prompt text.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs($count - 1)
$boldPara.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($count)
$r = $italicPara.Range
# Exclude the trailing paragraph mark so we only replace the run text.
$r.End = $r.End - 1
$r.Text = 'Please create a feature image fitting the game "Deadly 5" with the following specifications: - Cartoon style - Happy Maya warrior with glasses as the central character The image should convey the Wild West theme of the game and incorporate elements such as the dusty streets, saloon, cactus, and wooden barrel with handcuffs attached. The Maya warrior should be depicted as a confident and adventurous character, perhaps wielding a weapon or surrounded by other Wild West symbols. The image should also include the Deadly 5 logo and be visually appealing to attract players to the game.'
